# Trade #32 closed at 2026-02-17 08:28:24 - unknown UNKNOWN +0.000%
#
# Updates the aggregate statistics on the "Summary" and "Strategy Status"
# sheets to account for the newly-closed trade, and appends the trade's
# detail row (row 33) to both the "All Trades" and "MarketMaking" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B5").Value = -0.39    # Total P&L %
$summary.Range("B6").Value = 32       # Total Trades
$summary.Range("B9").Value = 21.88    # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet (MarketMaking row)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D4").Value = 32        # Trades
$status.Range("G4").Value = 21.88     # Win Rate %

# ---------------------------------------------------------------------
# Helper to append the new trade row (row 33) to a trade-log sheet
# ---------------------------------------------------------------------
function Add-TradeRow32($ws) {
    $row = 33

    $ws.Cells.Item($row, 1).Value = 32          # Trade #

    $ws.Cells.Item($row, 2).NumberFormat = "@"  # Date - keep as text (avoid auto date conversion)
    $ws.Cells.Item($row, 2).Value = "2026-02-17"

    $ws.Cells.Item($row, 3).Value = "08:28:18"  # Time - plain string, no auto-conversion

    $ws.Cells.Item($row, 4).Value = "MarketMaking"  # Strategy
    $ws.Cells.Item($row, 5).Value = "DOWN"          # Side
    $ws.Cells.Item($row, 6).Value = 0.02            # Entry Price
    $ws.Cells.Item($row, 7).Value = 0.021277        # Exit Price
    $ws.Cells.Item($row, 8).Value = "CLOSED"        # Status
    $ws.Cells.Item($row, 9).Value = 6.383           # P&L %
    $ws.Cells.Item($row, 10).Value = 0              # P&L $
    $ws.Cells.Item($row, 11).Value = 99.37          # Capital After
    $ws.Cells.Item($row, 12).Value = 0              # Entry Slippage (bps)
    $ws.Cells.Item($row, 13).Value = 0              # Exit Slippage (bps)
    $ws.Cells.Item($row, 14).Value = 0.6            # Confidence
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"  # Entry Reason
    $ws.Cells.Item($row, 16).Value = "early_exit"   # Exit Reason
    $ws.Cells.Item($row, 17).Value = 0.13           # Duration (min)
}

$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow32 $allTrades

$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow32 $marketMaking
